# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) across the profession sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2610.4443
$ws.Range("I32").Value = 9000
$ws.Range("J32").Value = 1811.75
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 1811.75
$ws.Range("M32").Value = -8674
$ws.Range("N32").Value = -2463.75

$ws.Range("H137").Value = 1324.8889
$ws.Range("I137").Value = 1224.5264
$ws.Range("J137").Value = 1563.25
$ws.Range("K137").Value = 3673.5792
$ws.Range("L137").Value = 4689.75
$ws.Range("M137").Value = -1123.5792
$ws.Range("N137").Value = -9789.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 230.5
$ws.Range("I5").Value = 145.25
$ws.Range("J5").Value = 401
$ws.Range("K5").Value = 145.25
$ws.Range("L5").Value = 401
$ws.Range("M5").Value = -33.25
$ws.Range("N5").Value = -625

$ws.Range("H53").Value = 10000
$ws.Range("I53").Value = 10000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 10000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -9318

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 230.5
$ws.Range("I4").Value = 145.25
$ws.Range("J4").Value = 401
$ws.Range("K4").Value = 145.25
$ws.Range("L4").Value = 401
$ws.Range("M4").Value = -30.25
$ws.Range("N4").Value = -631

$ws.Range("H38").Value = 22000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 22000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 22000
$ws.Range("N38").Value = -22832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1927.0333
$ws.Range("I31").Value = 1242.125
$ws.Range("J31").Value = 4666.6665
$ws.Range("K31").Value = 1242.125
$ws.Range("L31").Value = 4666.6665
$ws.Range("M31").Value = -947.125
$ws.Range("N31").Value = -5256.6665

$ws.Range("H34").Value = 1927.0333
$ws.Range("I34").Value = 1242.125
$ws.Range("J34").Value = 4666.6665
$ws.Range("K34").Value = 1242.125
$ws.Range("L34").Value = 4666.6665
$ws.Range("M34").Value = -1040.125
$ws.Range("N34").Value = -5070.6665

$ws.Range("H58").Value = 1278765.4
$ws.Range("I58").Value = 1611791.1
$ws.Range("J58").Value = 2166.6667
$ws.Range("K58").Value = 1611791.1
$ws.Range("L58").Value = 2166.6667
$ws.Range("M58").Value = -1611588.1
$ws.Range("N58").Value = -2572.6667

$ws.Range("H94").Value = 2010.1111
$ws.Range("I94").Value = 1603.1666
$ws.Range("J94").Value = 2213.5833
$ws.Range("K94").Value = 1603.1666
$ws.Range("L94").Value = 2213.5833
$ws.Range("M94").Value = -1152.1666
$ws.Range("N94").Value = -3115.5833

$ws.Range("H105").Value = 11890
$ws.Range("I105").Value = 13276.25
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 13276.25
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = -11529.25
$ws.Range("N105").Value = -4294

$ws.Range("H136").Value = 1278765.4
$ws.Range("I136").Value = 1611791.1
$ws.Range("J136").Value = 2166.6667
$ws.Range("K136").Value = 4835373.300000001
$ws.Range("L136").Value = 6500.000100000001
$ws.Range("M136").Value = -4832823.300000001
$ws.Range("N136").Value = -11600.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 529253
$ws.Range("I63").Value = 1051006
$ws.Range("J63").Value = 7500
$ws.Range("K63").Value = 3153018
$ws.Range("L63").Value = 22500
$ws.Range("M63").Value = -3152269
$ws.Range("N63").Value = -23998

$ws.Range("H64").Value = 4942.143
$ws.Range("I64").Value = 2756
$ws.Range("J64").Value = 5172.263
$ws.Range("K64").Value = 8268
$ws.Range("L64").Value = 15516.789
$ws.Range("M64").Value = -7998
$ws.Range("N64").Value = -16056.789

$ws.Range("H66").Value = 529253
$ws.Range("I66").Value = 1051006
$ws.Range("J66").Value = 7500
$ws.Range("K66").Value = 9459054
$ws.Range("L66").Value = 67500
$ws.Range("M66").Value = -9455310
$ws.Range("N66").Value = -74988

$ws.Range("H67").Value = 4942.143
$ws.Range("I67").Value = 2756
$ws.Range("J67").Value = 5172.263
$ws.Range("K67").Value = 8268
$ws.Range("L67").Value = 15516.789
$ws.Range("M67").Value = -7332
$ws.Range("N67").Value = -17388.789

$ws.Range("H68").Value = 976.2727
$ws.Range("I68").Value = 1020
$ws.Range("J68").Value = 951.2857
$ws.Range("K68").Value = 3060
$ws.Range("L68").Value = 2853.8571
$ws.Range("M68").Value = -2249
$ws.Range("N68").Value = -4475.8571

$ws.Range("H70").Value = 26506
$ws.Range("I70").Value = 26506
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 79518
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -79203

$ws.Range("H71").Value = 976.2727
$ws.Range("I71").Value = 1020
$ws.Range("J71").Value = 951.2857
$ws.Range("K71").Value = 9180
$ws.Range("L71").Value = 8561.5713
$ws.Range("M71").Value = -5124
$ws.Range("N71").Value = -16673.5713

$ws.Range("H73").Value = 26506
$ws.Range("I73").Value = 26506
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 79518
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -78426

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()

$ws.Range("H87").Value = 8121.973
$ws.Range("I87").Value = 1085.6666
$ws.Range("J87").Value = 9483.839
$ws.Range("K87").Value = 3256.9998
$ws.Range("L87").Value = 28451.517
$ws.Range("M87").Value = -2008.9998
$ws.Range("N87").Value = -30947.517

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()

$ws.Range("H90").Value = 8121.973
$ws.Range("I90").Value = 1085.6666
$ws.Range("J90").Value = 9483.839
$ws.Range("K90").Value = 9770.999400000001
$ws.Range("L90").Value = 85354.55100000001
$ws.Range("M90").Value = -3530.999400000001
$ws.Range("N90").Value = -97834.55100000001

$ws.Range("H98").Value = 1517
$ws.Range("I98").Value = 2479.8
$ws.Range("J98").Value = 714.6667
$ws.Range("K98").Value = 7439.400000000001
$ws.Range("L98").Value = 2144.0001
$ws.Range("M98").Value = -5941.400000000001
$ws.Range("N98").Value = -5140.0001

$ws.Range("H114").Value = 1282.5
$ws.Range("I114").Value = 650
$ws.Range("J114").Value = 1915
$ws.Range("K114").Value = 1950
$ws.Range("L114").Value = 5745
$ws.Range("M114").Value = 1304
$ws.Range("N114").Value = -12253

$ws.Range("H121").Value = 15932.849
$ws.Range("I121").Value = 746.44446
$ws.Range("J121").Value = 21627.75
$ws.Range("K121").Value = 2239.33338
$ws.Range("L121").Value = 64883.25
$ws.Range("M121").Value = -929.33338
$ws.Range("N121").Value = -67503.25

$ws.Range("H122").Value = 1135.5
$ws.Range("I122").Value = 996.6667
$ws.Range("J122").Value = 1218.8
$ws.Range("K122").Value = 8970.0003
$ws.Range("L122").Value = 10969.2
$ws.Range("M122").Value = -6520.0003
$ws.Range("N122").Value = -15869.2

$ws.Range("H129").Value = 4546598.5
$ws.Range("I129").Value = 915
$ws.Range("J129").Value = 5556750.5
$ws.Range("K129").Value = 2745
$ws.Range("L129").Value = 16670251.5
$ws.Range("M129").Value = 2255
$ws.Range("N129").Value = -16680251.5

$ws.Range("H131").Value = 3740.4695
$ws.Range("I131").Value = 11496.556
$ws.Range("J131").Value = 1995.35
$ws.Range("K131").Value = 34489.66800000001
$ws.Range("L131").Value = 5986.049999999999
$ws.Range("M131").Value = -29449.66800000001
$ws.Range("N131").Value = -16066.05

$ws.Range("H134").Value = 3458.4
$ws.Range("I134").Value = 1974.3889
$ws.Range("J134").Value = 5684.4165
$ws.Range("K134").Value = 5923.1667
$ws.Range("L134").Value = 17053.2495
$ws.Range("M134").Value = -853.1666999999998
$ws.Range("N134").Value = -27193.2495

$ws.Range("H136").Value = 5216.6
$ws.Range("I136").Value = 986.6667
$ws.Range("J136").Value = 8036.5557
$ws.Range("K136").Value = 2960.0001
$ws.Range("L136").Value = 24109.6671
$ws.Range("M136").Value = 2139.9999
$ws.Range("N136").Value = -34309.6671

$ws.Range("H139").Value = 1988.3334
$ws.Range("I139").Value = 1644.25
$ws.Range("J139").Value = 2971.4285
$ws.Range("K139").Value = 4932.75
$ws.Range("L139").Value = 8914.2855
$ws.Range("M139").Value = 207.25
$ws.Range("N139").Value = -19194.2855

$ws.Range("H140").Value = 5109.6665
$ws.Range("I140").Value = 1472.5
$ws.Range("J140").Value = 9959.223
$ws.Range("K140").Value = 4417.5
$ws.Range("L140").Value = 29877.669
$ws.Range("M140").Value = 762.5
$ws.Range("N140").Value = -40237.669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 214800
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 214800
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 214800
$ws.Range("N15").Value = -215376

$ws.Range("H52").Value = 500000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 500000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 500000
$ws.Range("N52").Value = -500518

$ws.Range("H55").Value = 6600
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 7000
$ws.Range("K55").Value = 5000
$ws.Range("L55").Value = 7000
$ws.Range("M55").Value = -4673
$ws.Range("N55").Value = -7654

$ws.Range("H81").Value = 214800
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 214800
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 214800
$ws.Range("N81").Value = -216796

$ws.Range("H84").Value = 214800
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 214800
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 644400
$ws.Range("N84").Value = -654384

$ws.Range("H131").Value = 31348.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 31348.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 31348.5
$ws.Range("N131").Value = -41428.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 20000
$ws.Range("I45").Value = 20000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 20000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -19593

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 39787.53
$ws.Range("I81").Value = 36088.83
$ws.Range("J81").Value = 61240
$ws.Range("K81").Value = 72177.66
$ws.Range("L81").Value = 122480
$ws.Range("M81").Value = -71116.66
$ws.Range("N81").Value = -124602

$ws.Range("H84").Value = 39787.53
$ws.Range("I84").Value = 36088.83
$ws.Range("J84").Value = 61240
$ws.Range("K84").Value = 360888.3
$ws.Range("L84").Value = 612400
$ws.Range("M84").Value = -355584.3
$ws.Range("N84").Value = -623008

$ws.Range("H123").Value = 35928.25
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 35928.25
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 35928.25
$ws.Range("N123").Value = -45728.25

$ws.Range("H125").Value = 46249
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 46249
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 46249
$ws.Range("N125").Value = -56089
